# SalesOrderConfirmationEmail.docx — "Syncing with version 27.0.36535.0"
#
# The Business Central report layout's embedded data-source schema
# (customXml part, namespace urn:microsoft-dynamics-nav/reports/
# Sales_Order_Confirmation_CZL/31187/) was regenerated by the NAV/BC
# Word-layout tooling:
#   1. A CRLF was inserted right after the XML declaration, putting the
#      root element <NavWordReportXmlPart> on its own line.
#   2. Two fields that are no longer emitted by the report data set were
#      dropped: DueDateFormat_SalesHeader and DueDate_SalesHeaderCaption.
#   3. The empty <Sales_Line/> element was re-serialized with a space
#      before its self-closing slash: <Sales_Line />.
#
# None of this is visible document body text (the report's field
# placeholders keep their existing bindings/captions in the document
# body) — it only affects the bound CustomXMLPart's raw schema markup,
# so we go straight at CustomXMLParts rather than Content.Find.

$d = $word.ActiveDocument

$targetNamespace = "urn:microsoft-dynamics-nav/reports/Sales_Order_Confirmation_CZL/31187/"
$targetPart = $null

# Real Word exposes CustomXMLParts.SelectByNamespace for exactly this
# ("find the NAV/BC report data part") — prefer it, and fall back to a
# manual scan in case a given host only supports enumeration.
try {
    $matches = $d.CustomXMLParts.SelectByNamespace($targetNamespace)
    if ($matches -ne $null -and $matches.Count -ge 1) {
        $targetPart = $matches.Item(1)
    }
} catch {
    $targetPart = $null
}

if ($targetPart -eq $null) {
    $parts = $d.CustomXMLParts
    for ($i = 1; $i -le $parts.Count; $i++) {
        $candidate = $parts.Item($i)
        if ($candidate.XML -ne $null -and $candidate.XML.Contains($targetNamespace)) {
            $targetPart = $candidate
            break
        }
    }
}

if ($targetPart -ne $null) {
    $xml = $targetPart.XML

    # 1) Break the xml declaration and the root element onto separate lines.
    $xml = $xml.Replace(
        "<?xml version=`"1.0`" encoding=`"utf-16`"?><NavWordReportXmlPart",
        "<?xml version=`"1.0`" encoding=`"utf-16`"?>`r`n<NavWordReportXmlPart")

    # 2) Drop the two retired fields (and their line terminators).
    $xml = $xml.Replace(
        "    <DueDateFormat_SalesHeader>DueDateFormat_SalesHeader</DueDateFormat_SalesHeader>`r`n    <DueDate_SalesHeaderCaption>DueDate_SalesHeaderCaption</DueDate_SalesHeaderCaption>`r`n",
        "")

    # 3) Re-pad the empty Sales_Line element's self-closing tag.
    $xml = $xml.Replace("<Sales_Line/>", "<Sales_Line />")

    $targetPart.XML = $xml
}
